# Sura - DataSource Emision Motor
# Adds 6 new rows (12-17) for the "Movilidad" item (gw environment /
# PQM042..PQM047 patente/motor/chasis triples), each with a hyperlink
# in column C pointing at the gw PolicyCenter environment.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$envDomain = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$envUrl    = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"

$rows = @(
    @{ Row = 12; F = 7166033480; Patente = "PQM042"; Motor = "ABC12PQM042"; Chasis = "ZAZ123PQM042" },
    @{ Row = 13; F = 7166033480; Patente = "PQM043"; Motor = "ABC12PQM043"; Chasis = "ZAZ123PQM043" },
    @{ Row = 14; F = 7166033480; Patente = "PQM044"; Motor = "ABC12PQM044"; Chasis = "ZAZ123PQM044" },
    @{ Row = 15; F = 1571314844; Patente = "PQM045"; Motor = "ABC12PQM045"; Chasis = "ZAZ123PQM045" },
    @{ Row = 16; F = 1571314844; Patente = "PQM046"; Motor = "ABC12PQM046"; Chasis = "ZAZ123PQM046" },
    @{ Row = 17; F = 1571314844; Patente = "PQM047"; Motor = "ABC12PQM047"; Chasis = "ZAZ123PQM047" }
)

foreach ($r in $rows) {
    $i = $r.Row

    $ws.Range("B$i").Value = "'$envDomain"
    $ws.Range("C$i").Value = $envUrl
    $ws.Hyperlinks.Add($ws.Range("C$i"), $envUrl)
    $ws.Range("D$i").Value = "su"
    $ws.Range("E$i").Value = "gw"
    $ws.Range("F$i").Value = $r.F
    $ws.Range("G$i").Value = "Motor"
    $ws.Range("H$i").Value = "Menos de 5 vehículos"
    $ws.Range("I$i").Value = "Anual"
    $ws.Range("J$i").Value = "'14/09/2020"
    $ws.Range("K$i").Value = "Cupón"
    $ws.Range("L$i").Value = "No"

    $ws.Range("S$i").Value = 2021
    $ws.Range("T$i").Value = "si"
    $ws.Range("U$i").Value = "CHEVROLET"
    $ws.Range("V$i").Value = "ONIX 1.2 L/19"
    $ws.Range("W$i").Value = 1487000
    $ws.Range("X$i").Value = "B - Resp. Civil-Robo/Incendio Total Daños Totales por Accidente"
    $ws.Range("Y$i").Value = $r.Patente
    $ws.Range("Z$i").Value = $r.Motor
    $ws.Range("AA$i").Value = $r.Chasis
    $ws.Range("AB$i").Value = "Sin Actividad"
    $ws.Range("AC$i").Value = "No"
}

$ws.Range("AB12:AC17").Select()
